$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 223 and 224 (columns B:AD), keep column A as-is ---
$r1 = 223
$r2 = 224
$cols = 2..30   # B=2 .. AD=30
foreach ($c in $cols) {
    $v1 = $ws.Cells.Item($r1, $c).Value2
    $v2 = $ws.Cells.Item($r2, $c).Value2
    $ws.Cells.Item($r1, $c).Value2 = $v2
    $ws.Cells.Item($r2, $c).Value2 = $v1
}

# --- Swap rows 230 and 232 (columns B:AD), keep column A as-is ---
$r3 = 230
$r4 = 232
foreach ($c in $cols) {
    $v3 = $ws.Cells.Item($r3, $c).Value2
    $v4 = $ws.Cells.Item($r4, $c).Value2
    $ws.Cells.Item($r3, $c).Value2 = $v4
    $ws.Cells.Item($r4, $c).Value2 = $v3
}
